$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 58.46
$ws.Range("C2").Value = 0.311
$ws.Range("D2").Value = 2007
$ws.Range("E2").Value = 57.851
$ws.Range("F2").Value = 59.069
$ws.Range("B3").Value = 64.71899999999999
$ws.Range("C3").Value = 0.307
$ws.Range("D3").Value = 2007
$ws.Range("E3").Value = 64.117
$ws.Range("F3").Value = 65.321
$ws.Range("B4").Value = 108.381
$ws.Range("C4").Value = 0.546
$ws.Range("D4").Value = 2007
$ws.Range("E4").Value = 107.311
$ws.Range("F4").Value = 109.451
$ws.Range("B5").Value = 112.815
$ws.Range("C5").Value = 0.539
$ws.Range("D5").Value = 2007
$ws.Range("E5").Value = 111.758
$ws.Range("F5").Value = 113.873
$ws.Range("B6").Value = 127.123
$ws.Range("C6").Value = 0.699
$ws.Range("D6").Value = 2007
$ws.Range("E6").Value = 125.753
$ws.Range("F6").Value = 128.494
$ws.Range("B7").Value = 136.58
$ws.Range("C7").Value = 0.6909999999999999
$ws.Range("D7").Value = 2007
$ws.Range("E7").Value = 135.225
$ws.Range("F7").Value = 137.935
$ws.Range("B8").Value = 94.19799999999999
$ws.Range("C8").Value = 0.76
$ws.Range("D8").Value = 2007
$ws.Range("E8").Value = 92.70699999999999
$ws.Range("F8").Value = 95.68899999999999
$ws.Range("B9").Value = 106.098
$ws.Range("C9").Value = 0.751
$ws.Range("D9").Value = 2007
$ws.Range("E9").Value = 104.624
$ws.Range("F9").Value = 107.571
$ws.Range("B10").Value = 18.134
$ws.Range("C10").Value = 0.269
$ws.Range("D10").Value = 2007
$ws.Range("E10").Value = 17.605
$ws.Range("F10").Value = 18.662
$ws.Range("B11").Value = 17.635
$ws.Range("C11").Value = 0.266
$ws.Range("D11").Value = 2007
$ws.Range("E11").Value = 17.112
$ws.Range("F11").Value = 18.157
$ws.Range("B12").Value = 18.228
$ws.Range("C12").Value = 0.139
$ws.Range("D12").Value = 2007
$ws.Range("E12").Value = 17.955
$ws.Range("F12").Value = 18.501
$ws.Range("B13").Value = 19.285
$ws.Range("C13").Value = 0.138
$ws.Range("D13").Value = 2007
$ws.Range("E13").Value = 19.015
$ws.Range("F13").Value = 19.554
$ws.Range("B14").Value = 42.833
$ws.Range("C14").Value = 0.227
$ws.Range("D14").Value = 2007
$ws.Range("E14").Value = 42.388
$ws.Range("F14").Value = 43.278
$ws.Range("B15").Value = 46.497
$ws.Range("C15").Value = 0.224
$ws.Range("D15").Value = 2007
$ws.Range("E15").Value = 46.058
$ws.Range("F15").Value = 46.937
$ws.Range("B16").Value = 111.82
$ws.Range("C16").Value = 0.479
$ws.Range("D16").Value = 2007
$ws.Range("E16").Value = 110.882
$ws.Range("F16").Value = 112.759
$ws.Range("B17").Value = 119.824
$ws.Range("C17").Value = 0.473
$ws.Range("D17").Value = 2007
$ws.Range("E17").Value = 118.896
$ws.Range("F17").Value = 120.752
$ws.Range("B18").Value = 71.706
$ws.Range("C18").Value = 0.548
$ws.Range("D18").Value = 2007
$ws.Range("E18").Value = 70.63200000000001
$ws.Range("F18").Value = 72.78100000000001
$ws.Range("B19").Value = 77.246
$ws.Range("C19").Value = 0.542
$ws.Range("D19").Value = 2007
$ws.Range("E19").Value = 76.184
$ws.Range("F19").Value = 78.30800000000001
$ws.Range("B20").Value = 147.151
$ws.Range("C20").Value = 0.667
$ws.Range("D20").Value = 2007
$ws.Range("E20").Value = 145.844
$ws.Range("F20").Value = 148.458
$ws.Range("B21").Value = 161.819
$ws.Range("C21").Value = 0.659
$ws.Range("D21").Value = 2007
$ws.Range("E21").Value = 160.527
$ws.Range("F21").Value = 163.111
$ws.Range("B22").Value = 276.408
$ws.Range("C22").Value = 0.703
$ws.Range("D22").Value = 2007
$ws.Range("E22").Value = 275.03
$ws.Range("F22").Value = 277.786
$ws.Range("B23").Value = 290.796
$ws.Range("C23").Value = 0.694
$ws.Range("D23").Value = 2007
$ws.Range("E23").Value = 289.434
$ws.Range("F23").Value = 292.157
$ws.Range("B24").Value = 142.857
$ws.Range("C24").Value = 0.347
$ws.Range("D24").Value = 2007
$ws.Range("E24").Value = 142.177
$ws.Range("F24").Value = 143.537
$ws.Range("B25").Value = 152.42
$ws.Range("C25").Value = 0.343
$ws.Range("D25").Value = 2007
$ws.Range("E25").Value = 151.748
$ws.Range("F25").Value = 153.092
